$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.135.65'
$ws.Range("D3").Value = '1.788.32'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.14'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +1.92%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.64'
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.19'
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.283'
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0663'
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '2.044.97'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.50'
$ws.Range("E14").Value = '  +13.68%  '
$ws.Range("D15").Value = '1.771.05'
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.634'
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '34.130.03'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.25'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.55'
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '254.44'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '0.0₃0744'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.52'
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.29'
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.61'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.06'
$ws.Range("E28").Value = '  +0.33%  '
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("E32").Value = '  +2.00%  '
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E34").Value = '  +2.70%  '
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("D36").Value = '1.452.48'
$ws.Range("E36").Value = '  -6.31%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.637'
$ws.Range("E38").Value = '  +3.89%  '
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +3.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '83.57'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.904'
$ws.Range("E43").Value = '  +0.48%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.85'
$ws.Range("E47").Value = '  +3.43%  '
$ws.Range("D48").Value = '1.944.94'
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.02'
$ws.Range("E49").Value = '  +9.13%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.20'
$ws.Range("E51").Value = '  -1.04%  '
